$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: "Pagos nuevo"
$ws.Cells.Item(34, 1).Value = 34
$ws.Cells.Item(34, 2).Value = 30
$ws.Cells.Item(34, 3).Value = "Pagos nuevo"
$ws.Cells.Item(34, 4).Value = "pagos/nuevo"
$ws.Cells.Item(34, 5).Value = "minimize"
$ws.Cells.Item(34, 6).Value = "oculto"
$ws.Cells.Item(34, 7).Value = "Digitador"
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = "Ruta para generar un nuevo pago"

# Row 35: "Pagos detalle"
$ws.Cells.Item(35, 1).Value = 35
$ws.Cells.Item(35, 2).Value = 30
$ws.Cells.Item(35, 3).Value = "Pagos detalle"
$ws.Cells.Item(35, 4).Value = "pagos/detalle"
$ws.Cells.Item(35, 5).Value = "minimize"
$ws.Cells.Item(35, 6).Value = "oculto"
$ws.Cells.Item(35, 7).Value = "Digitador"
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(35, 9).Value = "Ruta para el detalle de pagos"

# Update the view/selection to match the scrolled-down state after adding rows.
$ws.Range("A32").Select()
